$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("H112").Value = 2251.2307
$ws.Range("J112").Value = 2251.2307
$ws.Range("L112").Value = 6753.6921
$ws.Range("N112").Value = -8969.6921
$ws.Range("H129").Value = 1030.762
$ws.Range("J129").Value = 1144.4688
$ws.Range("L129").Value = 3433.4064
$ws.Range("N129").Value = -13433.4064
$ws.Range("H132").Value = 2095.074
$ws.Range("I132").Value = 1387.2821
$ws.Range("J132").Value = 3935.3333
$ws.Range("K132").Value = 4161.846299999999
$ws.Range("L132").Value = 11805.9999
$ws.Range("M132").Value = -1631.846299999999
$ws.Range("N132").Value = -16865.9999
$ws.Range("H138").Value = 2200896.5
$ws.Range("I138").Value = 9525517
$ws.Range("J138").Value = 3510.1572
$ws.Range("K138").Value = 28576551
$ws.Range("L138").Value = 10530.4716
$ws.Range("M138").Value = -28571411
$ws.Range("N138").Value = -20810.4716
$ws.Range("H141").Value = 21969
$ws.Range("I141").Value = 12316.111
$ws.Range("J141").Value = 36448.332
$ws.Range("K141").Value = 36948.333
$ws.Range("L141").Value = 109344.996
$ws.Range("M141").Value = -31768.333
$ws.Range("N141").Value = -119704.996

$ws = $wb.Worksheets.Item(2)
$ws.Range("H28").Value = 35996.707
$ws.Range("I28").Value = 36796.266
$ws.Range("K28").Value = 36796.266
$ws.Range("M28").Value = -36604.266
$ws.Range("H45").Value = 1148.1666
$ws.Range("I45").Value = 1047.25
$ws.Range("J45").Value = 1350
$ws.Range("K45").Value = 1047.25
$ws.Range("L45").Value = 1350
$ws.Range("M45").Value = -670.25
$ws.Range("N45").Value = -2104
$ws.Range("H58").Value = 40032.2
$ws.Range("J58").Value = 40032.2
$ws.Range("L58").Value = 40032.2
$ws.Range("N58").Value = -40892.2
$ws.Range("H63").Value = 3434.3157
$ws.Range("I63").Value = 3634.7693
$ws.Range("K63").Value = 3634.7693
$ws.Range("M63").Value = -2948.7693
$ws.Range("H66").Value = 3434.3157
$ws.Range("I66").Value = 3634.7693
$ws.Range("K66").Value = 18173.8465
$ws.Range("M66").Value = -14741.8465
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
[void]$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
[void]$ws.Range("N79").ClearContents()
$ws.Range("H80").Value = 36000
$ws.Range("J80").Value = 36000
$ws.Range("L80").Value = 36000
$ws.Range("N80").Value = -37996
$ws.Range("H82").Value = 27900
$ws.Range("J82").Value = 27900
$ws.Range("L82").Value = 27900
$ws.Range("N82").Value = -28622
$ws.Range("H83").Value = 36000
$ws.Range("J83").Value = 36000
$ws.Range("L83").Value = 108000
$ws.Range("N83").Value = -117984
$ws.Range("H85").Value = 27900
$ws.Range("J85").Value = 27900
$ws.Range("L85").Value = 27900
$ws.Range("N85").Value = -30396
$ws.Range("H86").Value = 30500
$ws.Range("J86").Value = 30500
$ws.Range("L86").Value = 30500
$ws.Range("N86").Value = -32872
$ws.Range("H89").Value = 30500
$ws.Range("J89").Value = 30500
$ws.Range("L89").Value = 91500
$ws.Range("N89").Value = -103356
$ws.Range("H93").Value = 29699.875
$ws.Range("J93").Value = 29699.875
$ws.Range("L93").Value = 29699.875
$ws.Range("N93").Value = -34691.875
$ws.Range("H94").Value = 21665
$ws.Range("J94").Value = 21665
$ws.Range("L94").Value = 21665
$ws.Range("N94").Value = -23467
$ws.Range("H96").Value = 275000
$ws.Range("J96").Value = 275000
$ws.Range("L96").Value = 275000
$ws.Range("N96").Value = -280492
$ws.Range("H97").Value = 655.6667
$ws.Range("I97").Value = 655.6667
$ws.Range("K97").Value = 655.6667
$ws.Range("M97").Value = -159.6667
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
[void]$ws.Range("N98").ClearContents()
$ws.Range("H99").Value = 35996.707
$ws.Range("I99").Value = 36796.266
$ws.Range("K99").Value = 36796.266
$ws.Range("M99").Value = -33801.266
$ws.Range("H132").Value = 1983.8064
$ws.Range("I132").Value = 1500
$ws.Range("J132").Value = 2863.4546
$ws.Range("K132").Value = 4500
$ws.Range("L132").Value = 8590.363799999999
$ws.Range("M132").Value = -1970
$ws.Range("N132").Value = -13650.3638

$ws = $wb.Worksheets.Item(4)
$ws.Range("I23").Value = 27375
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 27375
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -27135
[void]$ws.Range("N23").ClearContents()
$ws.Range("I27").Value = 27375
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 27375
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -27183
[void]$ws.Range("N27").ClearContents()
$ws.Range("H103").Value = 11444.833
$ws.Range("I103").Value = 2167.25
$ws.Range("J103").Value = 30000
$ws.Range("K103").Value = 2167.25
$ws.Range("L103").Value = 30000
$ws.Range("M103").Value = -995.25
$ws.Range("N103").Value = -32344
$ws.Range("H132").Value = 467511.66
$ws.Range("I132").Value = 712528.4399999999
$ws.Range("K132").Value = 2137585.32
$ws.Range("M132").Value = -2135055.32
$ws.Range("H134").Value = 1806.8235
$ws.Range("I134").Value = 1669
$ws.Range("J134").Value = 4012
$ws.Range("K134").Value = 5007
$ws.Range("L134").Value = 12036
$ws.Range("M134").Value = -2472
$ws.Range("N134").Value = -17106

$ws = $wb.Worksheets.Item(5)
$ws.Range("H107").Value = 803.4286
$ws.Range("I107").Value = 1080
$ws.Range("J107").Value = 434.66666
$ws.Range("K107").Value = 3240
$ws.Range("L107").Value = 1303.99998
$ws.Range("M107").Value = -1320
$ws.Range("N107").Value = -5143.999980000001
$ws.Range("H132").Value = 2358.5715
$ws.Range("I132").Value = 1101
$ws.Range("J132").Value = 5502.5
$ws.Range("K132").Value = 9909
$ws.Range("L132").Value = 49522.5
$ws.Range("M132").Value = -7379
$ws.Range("N132").Value = -54582.5

$ws = $wb.Worksheets.Item(6)
$ws.Range("H122").Value = 2619.8
$ws.Range("I122").Value = 2534.5652
$ws.Range("J122").Value = 3600
$ws.Range("K122").Value = 7603.6956
$ws.Range("L122").Value = 10800
$ws.Range("M122").Value = -5153.6956
$ws.Range("N122").Value = -15700

$ws = $wb.Worksheets.Item(7)
$ws.Range("H136").Value = 3265.8262
$ws.Range("I136").Value = 2565
$ws.Range("J136").Value = 5251.5
$ws.Range("K136").Value = 7695
$ws.Range("L136").Value = 15754.5
$ws.Range("M136").Value = -5145
$ws.Range("N136").Value = -20854.5
